$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3 content - "AutoAdminLogon" configuration entry added to the
# Configuration Version / Configuration Section / Notes / Examples table.
$ws.Range("A3").Value = 3.86
$ws.Range("B3").Value = '<AutoAdminLogon Enable="false" Password=""/> Configuration area'
$ws.Range("C3").Value = "The 'Password=""""' parameter is for the provision of the installer / setup account password you will be using should you want to enable auto logon after the machine reboots."
$ws.Range("D3").Value = '<AutoAdminLogon Enable="true" Password="DevPassword123"/>'

# Apply the same cell formatting used by row 2 (vertical-top for column A,
# vertical-top + wrap text for columns B:D).
$ws.Range("A3").VerticalAlignment = -4160
$ws.Range("B3:D3").WrapText = $true
$ws.Range("B3:D3").VerticalAlignment = -4160

# Match the row height used for the new wrapped row.
$ws.Rows.Item(3).RowHeight = 165

# Move the active selection to the newly added cell, as in the source workbook.
$ws.Range("D3").Select()
